$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.561.41"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "3.113.04"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.16"
$ws.Range("E5").Value = "  +8.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "623.96"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("E7").Value = "  -5.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.368"
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "3.107.53"
$ws.Range("E10").Value = "  -1.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.738"
$ws.Range("E11").Value = "  -1.75%  "
$ws.Range("E12").Value = "  -2.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.99"
$ws.Range("E14").Value = "  +2.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.49"
$ws.Range("E15").Value = "  -3.54%  "
$ws.Range("D16").Value = "90.250.87"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "3.661.51"
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("D18").Value = "3.097.44"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.87"
$ws.Range("E19").Value = "  +5.12%  "
$ws.Range("E20").Value = "  +5.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.08"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "438.29"
$ws.Range("E22").Value = "  -5.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.59"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.96"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.91"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.57"
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "89.01"
$ws.Range("E27").Value = "  -4.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.12"
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").Value = "3.245.07"
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.34"
$ws.Range("E31").Value = "  +2.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.161"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.196"
$ws.Range("E34").Value = "  +7.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.94"
$ws.Range("E35").Value = "  -4.20%  "
$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.85"
$ws.Range("E36").Value = "  +5.54%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.153"
$ws.Range("E37").Value = "  +7.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.30"
$ws.Range("E38").Value = "  +5.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "504.67"
$ws.Range("E39").Value = "  -3.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.91"
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.29"
$ws.Range("E41").Value = "  -1.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0893"
$ws.Range("E42").Value = "  +3.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.407"
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.44"
$ws.Range("E46").Value = "  +55.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.91"
$ws.Range("E47").Value = "  -3.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.693"
$ws.Range("E48").Value = "  +2.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "152.55"
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.82"
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("E51").Value = "  -0.56%  "
